# Apply the commit's change (regenerated IG metadata spreadsheet):
#  1. The "Title" row (A5) on the Metadata sheet previously had an empty
#     Value cell (B5). Set it to the Name value ("QualificationPAC"),
#     which already exists in the shared strings table.
#  2. The "Date" row (A8) Value cell (B8) is bumped to the new generation
#     timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B5").Value = "QualificationPAC"
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
